# Hjemme passive updated meanEMG legmaxROM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 1 header values (B1:E1) ---
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# --- Row 2: B2 and D2 updated; C2 and E2 removed entirely ---
$ws.Range("B2").Value = 1.5060282168894319
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 1.2385613962985595
$ws.Range("E2").ClearContents()

# --- Row 3: B3:E3 updated in place ---
$ws.Range("B3").Value = 0.9706601037273147
$ws.Range("C3").Value = -2.087878713081964
$ws.Range("D3").Value = 0.73228502515329164
$ws.Range("E3").Value = -2.561795311586474

# --- Update selection to match the new narrower data range ---
$ws.Range("B1:E3").Select()
